# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310   (left-hand / "before" comparison columns)
#   *_new -> *_FV2404   (right-hand / "after" comparison columns)
# Also turn the header row + data range into a proper Excel Table and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1, columns A:U) -----------------
$oldToNew = @{
    "Segmentname_old"          = "Segmentname_FV2310"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2310"
    "Segment_old"               = "Segment_FV2310"
    "Datenelement_old"         = "Datenelement_FV2310"
    "Segment ID_old"           = "Segment ID_FV2310"
    "Code_old"                  = "Code_FV2310"
    "Qualifier_old"             = "Qualifier_FV2310"
    "Beschreibung_old"         = "Beschreibung_FV2310"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2310"
    "Bedingung_old"             = "Bedingung_FV2310"
    "Segmentname_new"          = "Segmentname_FV2404"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2404"
    "Segment_new"               = "Segment_FV2404"
    "Datenelement_new"         = "Datenelement_FV2404"
    "Segment ID_new"           = "Segment ID_FV2404"
    "Code_new"                  = "Code_FV2404"
    "Qualifier_new"             = "Qualifier_FV2404"
    "Beschreibung_new"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2404"
    "Bedingung_new"             = "Bedingung_FV2404"
}

$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le 21; $c++) {
    $cell = $headerRange.Cells.Item(1, $c)
    $current = $cell.Value()
    if ($oldToNew.ContainsKey($current)) {
        $cell.Value = $oldToNew[$current]
    }
}

# --- 2. Turn the header+data range into an Excel Table ----------------
$dataRange = $ws.Range("A1:U78")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) ----------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
